$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Recommandations")
$ws2 = $wb.Worksheets.Item("Top_YTD")

# --- Sheet "Recommandations" updates ---
$ws1.Cells.Item(2,3).Value = 10
$ws1.Cells.Item(2,4).Value = 4182.41
$ws1.Cells.Item(2,5).Value = 104.85
$ws1.Cells.Item(3,3).Value = 5
$ws1.Cells.Item(3,4).Value = 3455
$ws1.Cells.Item(3,5).Value = 700
$ws1.Cells.Item(4,3).Value = 5
$ws1.Cells.Item(4,4).Value = 3405
$ws1.Cells.Item(4,5).Value = 680
$ws1.Cells.Item(5,3).Value = 5
$ws1.Cells.Item(5,4).Value = 3235.31
$ws1.Cells.Item(5,5).Value = 654.1900000000001
$ws1.Cells.Item(6,3).Value = 5
$ws1.Cells.Item(6,4).Value = 2134.29
$ws1.Cells.Item(6,5).Value = 429.04
$ws1.Cells.Item(7,3).Value = 5
$ws1.Cells.Item(7,4).Value = 1835.8
$ws1.Cells.Item(7,5).Value = 364.71
$ws1.Cells.Item(8,3).Value = 5
$ws1.Cells.Item(8,4).Value = 1751.08
$ws1.Cells.Item(8,5).Value = 366.4
$ws1.Cells.Item(9,3).Value = 5
$ws1.Cells.Item(9,4).Value = 702.33
$ws1.Cells.Item(9,5).Value = 141.75
$ws1.Cells.Item(10,1).Value = 'BRVM-PRESTIGE'
$ws1.Cells.Item(10,3).Value = 5
$ws1.Cells.Item(10,4).Value = 684.4
$ws1.Cells.Item(10,5).Value = 137.03
$ws1.Cells.Item(11,1).Value = 'BRVM - INDUSTRIELS'
$ws1.Cells.Item(11,3).Value = 5
$ws1.Cells.Item(11,4).Value = 676.9299999999999
$ws1.Cells.Item(11,5).Value = 118.39
$ws1.Cells.Item(12,3).Value = 5
$ws1.Cells.Item(12,4).Value = 673.13
$ws1.Cells.Item(12,5).Value = 135.66
$ws1.Cells.Item(13,3).Value = 5
$ws1.Cells.Item(13,4).Value = 661.53
$ws1.Cells.Item(13,5).Value = 133.32
$ws1.Cells.Item(14,3).Value = 5
$ws1.Cells.Item(14,4).Value = 538.01
$ws1.Cells.Item(14,5).Value = 107.77
$ws1.Cells.Item(15,3).Value = 5
$ws1.Cells.Item(15,4).Value = 475.46
$ws1.Cells.Item(15,5).Value = 94.64
$ws1.Cells.Item(25,3).Value = 1
$ws1.Cells.Item(25,4).Value = 26.42
$ws1.Cells.Item(25,5).Value = -3.33
$ws1.Cells.Item(27,3).Value = 1
$ws1.Cells.Item(27,4).Value = 15.36
$ws1.Cells.Item(27,5).Value = -3.31
$ws1.Cells.Item(28,1).Value = 'SETAO CI (STAC)'
$ws1.Cells.Item(28,2).Value = 3
$ws1.Cells.Item(28,4).Value = 10.74
$ws1.Cells.Item(28,5).Value = 5.42
$ws1.Cells.Item(28,6).Value = '🟢 Achat'
$ws1.Cells.Item(28,7).Value = '✅ Renforcer'
$ws1.Cells.Item(29,1).Value = 'SAPH CI (SPHC)'
$ws1.Cells.Item(29,2).Value = 2
$ws1.Cells.Item(29,3).Value = 1
$ws1.Cells.Item(29,4).Value = 8.949999999999999
$ws1.Cells.Item(29,5).Value = 7.48
$ws1.Cells.Item(29,7).Value = '👀 À surveiller'
$ws1.Cells.Item(30,1).Value = 'SOGB CI (SOGC)'
$ws1.Cells.Item(30,2).Value = 1
$ws1.Cells.Item(30,3).Value = 0
$ws1.Cells.Item(30,4).Value = 6.97
$ws1.Cells.Item(30,5).Value = 6.97
$ws1.Cells.Item(30,7).Value = '➖ Neutre'
$ws1.Cells.Item(33,1).Value = 'SODE CI (SDCC)'
$ws1.Cells.Item(33,3).Value = 1
$ws1.Cells.Item(33,4).Value = 2.9
$ws1.Cells.Item(33,5).Value = 4.35
$ws1.Cells.Item(33,7).Value = '👀 À surveiller'
$ws1.Cells.Item(34,1).Value = 'BICI CI (BICC)'
$ws1.Cells.Item(34,3).Value = 0
$ws1.Cells.Item(34,4).Value = 1.96
$ws1.Cells.Item(34,5).Value = 1.96
$ws1.Cells.Item(34,7).Value = '➖ Neutre'
$ws1.Cells.Item(35,3).Value = 4
$ws1.Cells.Item(36,1).Value = 'ECOBANK COTE D''''IVOIRE (ECOC)'
$ws1.Cells.Item(36,4).Value = -1.5
$ws1.Cells.Item(36,5).Value = -1.5
$ws1.Cells.Item(37,1).Value = 'NEI-CEDA CI (NEIC)'
$ws1.Cells.Item(37,2).Value = 1
$ws1.Cells.Item(37,4).Value = -2.02
$ws1.Cells.Item(37,5).Value = 3.82
$ws1.Cells.Item(37,7).Value = '👀 À surveiller'
$ws1.Cells.Item(38,1).Value = 'ORANGE COTE D''IVOIRE (ORAC)'
$ws1.Cells.Item(38,2).Value = 0
$ws1.Cells.Item(38,4).Value = -2.03
$ws1.Cells.Item(38,5).Value = -2.03
$ws1.Cells.Item(38,7).Value = '➖ Neutre'
$ws1.Cells.Item(39,1).Value = 'UNILEVER CI (UNLC)'
$ws1.Cells.Item(39,2).Value = 1
$ws1.Cells.Item(39,4).Value = -2.23
$ws1.Cells.Item(39,5).Value = 5.26
$ws1.Cells.Item(39,7).Value = '👀 À surveiller'
$ws1.Cells.Item(46,1).Value = 'SICOR CI (SICC)'
$ws1.Cells.Item(46,2).Value = 1
$ws1.Cells.Item(46,3).Value = 2
$ws1.Cells.Item(46,4).Value = -6.47
$ws1.Cells.Item(46,5).Value = -6.91
$ws1.Cells.Item(46,7).Value = '👀 À surveiller'
$ws1.Cells.Item(47,1).Value = 'TRACTAFRIC MOTORS CI (PRSC)'
$ws1.Cells.Item(47,3).Value = 1
$ws1.Cells.Item(47,4).Value = -7.42
$ws1.Cells.Item(47,5).Value = -7.42
$ws1.Cells.Item(48,1).Value = 'BERNABE CI (BNBC)'
$ws1.Cells.Item(48,3).Value = 1
$ws1.Cells.Item(48,4).Value = -7.46
$ws1.Cells.Item(48,5).Value = -7.46
$ws1.Cells.Item(48,6).Value = '🟡 Observer'
$ws1.Cells.Item(48,7).Value = '➖ Neutre'
$ws1.Cells.Item(49,1).Value = 'BANK OF AFRICA BF (BOABF)'
$ws1.Cells.Item(49,2).Value = 1
$ws1.Cells.Item(49,3).Value = 2
$ws1.Cells.Item(49,4).Value = -7.48
$ws1.Cells.Item(49,5).Value = 4.85
$ws1.Cells.Item(49,6).Value = '🟡 Observer'
$ws1.Cells.Item(49,7).Value = '👀 À surveiller'
$ws1.Cells.Item(50,1).Value = 'FILTISAC CI (FTSC)'
$ws1.Cells.Item(50,2).Value = 1
$ws1.Cells.Item(50,3).Value = 3
$ws1.Cells.Item(50,4).Value = -12.21
$ws1.Cells.Item(50,5).Value = 7.37
$ws1.Cells.Item(50,6).Value = '🔴 Vente'
$ws1.Cells.Item(50,7).Value = '⚠️ Risque de décrochage'

# --- Sheet "Top_YTD" updates ---
$ws2.Cells.Item(2,2).Value = 149335173.34
$ws2.Cells.Item(3,2).Value = 3096060.89
$ws2.Cells.Item(4,2).Value = 2904741.51
$ws2.Cells.Item(5,2).Value = 2326605.9
$ws2.Cells.Item(6,2).Value = 405649.75
$ws2.Cells.Item(7,2).Value = 222392.98
$ws2.Cells.Item(8,2).Value = 184610.07
$ws2.Cells.Item(9,2).Value = 7935.48
$ws2.Cells.Item(10,1).Value = 'BRVM-PRESTIGE'
$ws2.Cells.Item(10,2).Value = 7358.31
$ws2.Cells.Item(11,1).Value = 'BRVM - INDUSTRIELS'
$ws2.Cells.Item(11,2).Value = 7098.57
